$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.803.77"
$ws.Range("E2").Value = "  +4.34%  "
$ws.Range("D3").Value = "3.266.50"
$ws.Range("E3").Value = "  +4.26%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'580.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.04%  "
$ws.Range("D6").Value = "'182.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.25%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.602"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.29%  "
$ws.Range("D9").Value = "3.266.52"
$ws.Range("E9").Value = "  +4.26%  "
$ws.Range("D10").Value = "'0.134"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.90%  "
$ws.Range("D11").Value = "'6.76"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.57%  "
$ws.Range("E12").Value = "  +7.33%  "
$ws.Range("D13").Value = "3.836.72"
$ws.Range("E13").Value = "  +4.35%  "
$ws.Range("E14").Value = "  +1.21%  "
$ws.Range("D15").Value = "'28.62"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.38%  "
$ws.Range("D16").Value = "67.803.69"
$ws.Range("E16").Value = "  +4.53%  "
$ws.Range("E17").Value = "  +5.02%  "
$ws.Range("D18").Value = "3.268.32"
$ws.Range("E18").Value = "  +4.30%  "
$ws.Range("E19").Value = "  +3.19%  "
$ws.Range("D20").Value = "'13.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.07%  "
$ws.Range("D21").Value = "'376.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.20%  "
$ws.Range("D22").Value = "'7.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.32%  "
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").Value = "'71.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.66%  "
$ws.Range("E25").Value = "  +4.37%  "
$ws.Range("E26").Value = "  +5.53%  "
$ws.Range("D27").Value = "'9.65"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("E28").Value = "  +3.51%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("E30").Value = "  +4.18%  "
$ws.Range("D31").Value = "'5.72"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.49%  "
$ws.Range("D32").Value = "'22.79"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.96%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("E34").Value = "  +8.06%  "
$ws.Range("D35").Value = "'6.94"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.03%  "
$ws.Range("E36").Value = "  +6.48%  "
$ws.Range("D37").Value = "'163.38"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.10%  "
$ws.Range("D38").Value = "'0.852"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.21%  "
$ws.Range("E39").Value = "  +5.78%  "
$ws.Range("E40").Value = "  +13.15%  "
$ws.Range("D41").Value = "'4.68"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +13.14%  "
$ws.Range("E42").Value = "  +2.86%  "
$ws.Range("D44").Value = "'354.74"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.90%  "
$ws.Range("D45").Value = "2.709.56"
$ws.Range("E45").Value = "  +2.68%  "
$ws.Range("D46").Value = "'25.54"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.31%  "
$ws.Range("D47").Value = "'40.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.99%  "
$ws.Range("E48").Value = "  +5.33%  "
$ws.Range("D49").Value = "'0.0281"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.18%  "
$ws.Range("E50").Value = "  +7.80%  "
$ws.Range("E51").Value = "  +0.80%  "
